$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the second scan entry (student 211267 / Manual) plus the trailing blank row,
# leaving only the header row and the single remaining log entry.
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()

# The surviving Log Time entry is refreshed to a plain text time stamp instead of
# the previous time-formatted serial number.
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = "12:54:41"
